# The commit adds a header value "Ou ID" to cell A1 of Sheet1
# (this pulls in a new shared-strings table) and leaves the
# active selection on A2, ready for further data entry below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Ou ID"
$ws.Range("A2").Select()
